$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B2: 10 -> 6
$ws.Range("B2").Value = 6

# Copy the formatting of A2 (style index 1: bordered, bold, centered) so it
# can be applied to the newly inserted A-column cells below.
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in new / shifted data rows
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1
